$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.552930018879517
$ws.Range("D2").Value = 8.149557194612376
$ws.Range("E2").Value = 13.25732633479951
$ws.Range("F2").Value = 37.79056479970345
$ws.Range("G2").Value = 43.10513593011221
$ws.Range("H2").Value = 17.72054311372988
$ws.Range("I2").Value = 21.35731689086278
$ws.Range("J2").Value = 10.10502759016736
$ws.Range("K2").Value = 14.85530324584142
$ws.Range("B3").Value = 5.429268336183933
$ws.Range("D3").Value = 8.108118459319751
$ws.Range("E3").Value = 13.20573248395244
$ws.Range("F3").Value = 37.78938057587226
$ws.Range("G3").Value = 43.05851015898448
$ws.Range("H3").Value = 17.76916508207933
$ws.Range("I3").Value = 21.4821306980641
$ws.Range("J3").Value = 10.11362009541285
$ws.Range("K3").Value = 14.48804979010264
$ws.Range("B4").Value = 5.35258890092231
$ws.Range("D4").Value = 8.083713616122761
$ws.Range("E4").Value = 13.17670741048135
$ws.Range("F4").Value = 37.79971956108497
$ws.Range("G4").Value = 43.04552441857165
$ws.Range("H4").Value = 17.80292537710034
$ws.Range("I4").Value = 21.5626802002113
$ws.Range("J4").Value = 10.12063183617648
$ws.Range("K4").Value = 14.26053489134801
$ws.Range("B5").Value = 5.321196299668718
$ws.Range("D5").Value = 8.074036589514336
$ws.Range("E5").Value = 13.16555497847319
$ws.Range("F5").Value = 37.80670982419447
$ws.Range("G5").Value = 43.04415922650757
$ws.Range("H5").Value = 17.81766219055466
$ws.Range("I5").Value = 21.59649105473644
$ws.Range("J5").Value = 10.1239255354741
$ws.Range("K5").Value = 14.16745120350679
$ws.Range("B6").Value = 5.31597612642778
$ws.Range("D6").Value = 8.072446118580869
$ws.Range("E6").Value = 13.16374414955715
$ws.Range("F6").Value = 37.80803807487035
$ws.Range("G6").Value = 43.04416940858334
$ws.Range("H6").Value = 17.82016827864129
$ws.Range("I6").Value = 21.60216495640679
$ws.Range("J6").Value = 10.1244988000399
$ws.Range("K6").Value = 14.151976568847
$ws.Range("B7").Value = 5.352166057936121
$ws.Range("D7").Value = 8.083582013701161
$ws.Range("E7").Value = 13.17655425946845
$ws.Range("F7").Value = 37.79980259912197
$ws.Range("G7").Value = 43.04549012091873
$ws.Range("H7").Value = 17.80312016225514
$ws.Range("I7").Value = 21.5631321888603
$ws.Range("J7").Value = 10.120674489648
$ws.Range("K7").Value = 14.25928084120442
$ws.Range("B8").Value = 5.510472918586279
$ws.Range("D8").Value = 8.135057803236011
$ws.Range("E8").Value = 13.23899115261819
$ws.Range("F8").Value = 37.78785692745659
$ws.Range("G8").Value = 43.08580990110261
$ws.Range("H8").Value = 17.7364954488279
$ws.Range("I8").Value = 21.3995418346211
$ws.Range("J8").Value = 10.10762993947368
$ws.Range("K8").Value = 14.72917458271818
$ws.Range("B9").Value = 5.813131430572724
$ws.Range("D9").Value = 8.243924173531919
$ws.Range("E9").Value = 13.38209004825429
$ws.Range("F9").Value = 37.85239893187715
$ws.Range("G9").Value = 43.28910006660438
$ws.Range("H9").Value = 17.63696967150023
$ws.Range("I9").Value = 21.10968858541722
$ws.Range("J9").Value = 10.09582651939031
$ws.Range("K9").Value = 15.6288564345665
$ws.Range("B10").Value = 6.0284549801062
$ws.Range("D10").Value = 8.328292370098843
$ws.Range("E10").Value = 13.49924051263121
$ws.Range("F10").Value = 37.95353841823655
$ws.Range("G10").Value = 43.5140397143682
$ws.Range("H10").Value = 17.58299535948141
$ws.Range("I10").Value = 20.91545304880825
$ws.Range("J10").Value = 10.09555263202793
$ws.Range("K10").Value = 16.26938719657515
$ws.Range("B11").Value = 6.124442598580352
$ws.Range("D11").Value = 8.367525447167949
$ws.Range("E11").Value = 13.55500621858627
$ws.Range("F11").Value = 38.01118908195101
$ws.Range("G11").Value = 43.63265573710759
$ws.Range("H11").Value = 17.56263279012043
$ws.Range("I11").Value = 20.83112423414241
$ws.Range("J11").Value = 10.09724897663059
$ws.Range("K11").Value = 16.55497511731584
$ws.Range("B12").Value = 6.160474783206904
$ws.Range("D12").Value = 8.382495940104961
$ws.Range("E12").Value = 13.57646622736699
$ws.Range("F12").Value = 38.03468802976667
$ws.Range("G12").Value = 43.67989695531821
$ws.Range("H12").Value = 17.55552715076233
$ws.Range("I12").Value = 20.79976838798393
$ws.Range("J12").Value = 10.09815270496816
$ws.Range("K12").Value = 16.66218456240085
$ws.Range("B13").Value = 6.152729184198655
$ws.Range("D13").Value = 8.379266860740502
$ws.Range("E13").Value = 13.57182939861493
$ws.Range("F13").Value = 38.02955305331317
$ws.Range("G13").Value = 43.66961970060873
$ws.Range("H13").Value = 17.55703052000727
$ws.Range("I13").Value = 20.80649576853126
$ws.Range("J13").Value = 10.09794645611143
$ws.Range("K13").Value = 16.63913830020183
$ws.Range("B14").Value = 6.127413510811626
$ws.Range("D14").Value = 8.368754850272561
$ws.Range("E14").Value = 13.55676495333292
$ws.Range("F14").Value = 38.01308897090595
$ws.Range("G14").Value = 43.63649584831057
$ws.Range("H14").Value = 17.56203606093777
$ws.Range("I14").Value = 20.82853300337432
$ws.Range("J14").Value = 10.09731809215939
$ws.Range("K14").Value = 16.56381463481776
$ws.Range("B15").Value = 6.111864785142473
$ws.Range("D15").Value = 8.36233048875622
$ws.Range("E15").Value = 13.54758178795442
$ws.Range("F15").Value = 38.00322123191124
$ws.Range("G15").Value = 43.61650854583734
$ws.Range("H15").Value = 17.56518099199902
$ws.Range("I15").Value = 20.84210661837088
$ws.Range("J15").Value = 10.09696722096268
$ws.Range("K15").Value = 16.51755182601697
$ws.Range("B16").Value = 6.022139493383339
$ws.Range("D16").Value = 8.325744809820389
$ws.Range("E16").Value = 13.49564477195848
$ws.Range("F16").Value = 37.95000459742543
$ws.Range("G16").Value = 43.50661418862744
$ws.Range("H16").Value = 17.58441065496405
$ws.Range("I16").Value = 20.92104506809713
$ws.Range("J16").Value = 10.09547837513882
$ws.Range("K16").Value = 16.25059774871679
$ws.Range("B17").Value = 5.966568025791917
$ws.Range("D17").Value = 8.303512976967991
$ws.Range("E17").Value = 13.46440757913223
$ws.Range("F17").Value = 37.92033597669248
$ws.Range("G17").Value = 43.44335801591705
$ws.Range("H17").Value = 17.59728266528927
$ws.Range("I17").Value = 20.97050207624563
$ws.Range("J17").Value = 10.09503109744189
$ws.Range("K17").Value = 16.08527008478584
$ws.Range("B18").Value = 5.934421775655519
$ws.Range("D18").Value = 8.290806704306329
$ws.Range("E18").Value = 13.44667438994836
$ws.Range("F18").Value = 37.90436766655741
$ws.Range("G18").Value = 43.40850941780053
$ws.Range("H18").Value = 17.60508060008631
$ws.Range("I18").Value = 20.99932786110949
$ws.Range("J18").Value = 10.09494526831799
$ws.Range("K18").Value = 15.98963839826362
$ws.Range("B19").Value = 5.923507246131175
$ws.Range("D19").Value = 8.286518741818433
$ws.Range("E19").Value = 13.44071074494411
$ws.Range("F19").Value = 37.89914951235563
$ws.Range("G19").Value = 43.39697435504051
$ws.Range("H19").Value = 17.60778848468058
$ws.Range("I19").Value = 21.00915299962748
$ws.Range("J19").Value = 10.09494566451849
$ws.Range("K19").Value = 15.95716976990753
$ws.Range("B20").Value = 5.97250289718358
$ws.Range("D20").Value = 8.305871286141727
$ws.Range("E20").Value = 13.46770874482238
$ws.Range("F20").Value = 37.92338082749906
$ws.Range("G20").Value = 43.44993303046454
$ws.Range("H20").Value = 17.59587159041271
$ws.Range("I20").Value = 20.96519803822289
$ws.Range("J20").Value = 10.0950609713675
$ws.Range("K20").Value = 16.10292611866285
$ws.Range("B21").Value = 6.134858177492679
$ws.Range("D21").Value = 8.371839466904952
$ws.Range("E21").Value = 13.56118054987299
$ws.Range("F21").Value = 38.01787966531521
$ws.Range("G21").Value = 43.64616222670414
$ws.Range("H21").Value = 17.56054936663097
$ws.Range("I21").Value = 20.82204447070198
$ws.Range("J21").Value = 10.0974955692389
$ws.Range("K21").Value = 16.58596522088396
$ws.Range("B22").Value = 6.23910951113425
$ws.Range("D22").Value = 8.415612405820111
$ws.Range("E22").Value = 13.62426117867986
$ws.Range("F22").Value = 38.08935742058014
$ws.Range("G22").Value = 43.78794309362857
$ws.Range("H22").Value = 17.54099304806429
$ws.Range("I22").Value = 20.73185147214712
$ws.Range("J22").Value = 10.10060979905718
$ws.Range("K22").Value = 16.89615843501431
$ws.Range("B23").Value = 6.183649070952803
$ws.Range("D23").Value = 8.392192675348369
$ws.Range("E23").Value = 13.59041606095032
$ws.Range("F23").Value = 38.05032188735918
$ws.Range("G23").Value = 43.71104092326956
$ws.Range("H23").Value = 17.55110691889725
$ws.Range("I23").Value = 20.77968177180113
$ws.Range("J23").Value = 10.09880851106497
$ws.Range("K23").Value = 16.7311378301814
$ws.Range("B24").Value = 5.969820354328211
$ws.Range("D24").Value = 8.304804859958349
$ws.Range("E24").Value = 13.46621558449121
$ws.Range("F24").Value = 37.9220008591951
$ws.Range("G24").Value = 43.44695573593714
$ws.Range("H24").Value = 17.59650829892038
$ws.Range("I24").Value = 20.96759477171392
$ws.Range("J24").Value = 10.09504693170549
$ws.Range("K24").Value = 16.09494563248269
$ws.Range("B25").Value = 5.73233311627019
$ws.Range("D25").Value = 8.213669708196267
$ws.Range("E25").Value = 13.34122307107216
$ws.Range("F25").Value = 37.82550284960551
$ws.Range("G25").Value = 43.22080287479503
$ws.Range("H25").Value = 17.66054348114377
$ws.Range("I25").Value = 21.18480376476461
$ws.Range("J25").Value = 10.09754431064214
$ws.Range("K25").Value = 15.3885727167179
